# Update scripts with new TPM values (Pgf-Flt1 ligand-receptor pairs).
# Adds the "Resolving-Mac" cluster (rows 5,9,13 gain it as target; new rows
# 14-17 add it as a sending cluster) and refreshes all computed NATMI metrics
# for rows 2-13 to match the new TPM-based recalculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Pgf"
$ws.Cells.Item(2, 3).Value = "Flt1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 3.529548333333333
$ws.Cells.Item(2, 8).Value = 10.588645
$ws.Cells.Item(2, 9).Value = 0.6379922832219154
$ws.Cells.Item(2, 10).Value = 0.6379922832219154
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 136.544502
$ws.Cells.Item(2, 14).Value = 409.633506
$ws.Cells.Item(2, 15).Value = 0.9681180443787725
$ws.Cells.Item(2, 16).Value = 0.9681180443787725
$ws.Cells.Item(2, 17).Value = 481.94041945993
$ws.Cells.Item(2, 18).Value = 4337.46377513937
$ws.Cells.Item(2, 19).Value = 0.6176518415615486
$ws.Cells.Item(2, 20).Value = 0.6176518415615486

# Row 3: ECs -> FAPs
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Pgf"
$ws.Cells.Item(3, 3).Value = "Flt1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 3.529548333333333
$ws.Cells.Item(3, 8).Value = 10.588645
$ws.Cells.Item(3, 9).Value = 0.6379922832219154
$ws.Cells.Item(3, 10).Value = 0.6379922832219154
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.3890956666666667
$ws.Cells.Item(3, 14).Value = 1.167287
$ws.Cells.Item(3, 15).Value = 0.002758738216274633
$ws.Cells.Item(3, 16).Value = 0.002758738216274633
$ws.Cells.Item(3, 17).Value = 1.373331961790555
$ws.Cells.Item(3, 18).Value = 12.359987656115
$ws.Cells.Item(3, 19).Value = 0.001760053693412607
$ws.Cells.Item(3, 20).Value = 0.001760053693412607

# Row 4: ECs -> MuSCs
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Pgf"
$ws.Cells.Item(4, 3).Value = "Flt1"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 3.529548333333333
$ws.Cells.Item(4, 8).Value = 10.588645
$ws.Cells.Item(4, 9).Value = 0.6379922832219154
$ws.Cells.Item(4, 10).Value = 0.6379922832219154
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 4.025396333333334
$ws.Cells.Item(4, 14).Value = 12.076189
$ws.Cells.Item(4, 15).Value = 0.0285405766544606
$ws.Cells.Item(4, 16).Value = 0.02854057665446059
$ws.Cells.Item(4, 17).Value = 14.20783091932278
$ws.Cells.Item(4, 18).Value = 127.870478273905
$ws.Cells.Item(4, 19).Value = 0.01820866766424941
$ws.Cells.Item(4, 20).Value = 0.01820866766424941

# Row 5: ECs -> Resolving-Mac
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Pgf"
$ws.Cells.Item(5, 3).Value = "Flt1"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 3.529548333333333
$ws.Cells.Item(5, 8).Value = 10.588645
$ws.Cells.Item(5, 9).Value = 0.6379922832219154
$ws.Cells.Item(5, 10).Value = 0.6379922832219154
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.08217633333333334
$ws.Cells.Item(5, 14).Value = 0.246529
$ws.Cells.Item(5, 15).Value = 0.0005826407504923545
$ws.Cells.Item(5, 16).Value = 0.0005826407504923544
$ws.Cells.Item(5, 17).Value = 0.2900453403561111
$ws.Cells.Item(5, 18).Value = 2.610408063205
$ws.Cells.Item(5, 19).Value = 0.0003717203027047476
$ws.Cells.Item(5, 20).Value = 0.0003717203027047475

# Row 6: FAPs -> ECs
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Pgf"
$ws.Cells.Item(6, 3).Value = "Flt1"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 0.6666666666666666
$ws.Cells.Item(6, 7).Value = 0.6384863333333334
$ws.Cells.Item(6, 8).Value = 1.915459
$ws.Cells.Item(6, 9).Value = 0.1154111844176443
$ws.Cells.Item(6, 10).Value = 0.1154111844176443
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 136.544502
$ws.Cells.Item(6, 14).Value = 409.633506
$ws.Cells.Item(6, 15).Value = 0.9681180443787725
$ws.Cells.Item(6, 16).Value = 0.9681180443787725
$ws.Cells.Item(6, 17).Value = 87.181798418806
$ws.Cells.Item(6, 18).Value = 784.636185769254
$ws.Cells.Item(6, 19).Value = 0.1117316501578476
$ws.Cells.Item(6, 20).Value = 0.1117316501578476

# Row 7: FAPs -> FAPs
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Pgf"
$ws.Cells.Item(7, 3).Value = "Flt1"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 0.6666666666666666
$ws.Cells.Item(7, 7).Value = 0.6384863333333334
$ws.Cells.Item(7, 8).Value = 1.915459
$ws.Cells.Item(7, 9).Value = 0.1154111844176443
$ws.Cells.Item(7, 10).Value = 0.1154111844176443
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.3890956666666667
$ws.Cells.Item(7, 14).Value = 1.167287
$ws.Cells.Item(7, 15).Value = 0.002758738216274633
$ws.Cells.Item(7, 16).Value = 0.002758738216274633
$ws.Cells.Item(7, 17).Value = 0.2484322655258889
$ws.Cells.Item(7, 18).Value = 2.235890389733
$ws.Cells.Item(7, 19).Value = 0.0003183892450384747
$ws.Cells.Item(7, 20).Value = 0.0003183892450384747

# Row 8: FAPs -> MuSCs
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Pgf"
$ws.Cells.Item(8, 3).Value = "Flt1"
$ws.Cells.Item(8, 4).Value = "MuSCs"
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 0.6666666666666666
$ws.Cells.Item(8, 7).Value = 0.6384863333333334
$ws.Cells.Item(8, 8).Value = 1.915459
$ws.Cells.Item(8, 9).Value = 0.1154111844176443
$ws.Cells.Item(8, 10).Value = 0.1154111844176443
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 4.025396333333334
$ws.Cells.Item(8, 14).Value = 12.076189
$ws.Cells.Item(8, 15).Value = 0.0285405766544606
$ws.Cells.Item(8, 16).Value = 0.02854057665446059
$ws.Cells.Item(8, 17).Value = 2.570160545083445
$ws.Cells.Item(8, 18).Value = 23.131444905751
$ws.Cells.Item(8, 19).Value = 0.003293901755653865
$ws.Cells.Item(8, 20).Value = 0.003293901755653864

# Row 9: FAPs -> Resolving-Mac
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Pgf"
$ws.Cells.Item(9, 3).Value = "Flt1"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 0.6666666666666666
$ws.Cells.Item(9, 7).Value = 0.6384863333333334
$ws.Cells.Item(9, 8).Value = 1.915459
$ws.Cells.Item(9, 9).Value = 0.1154111844176443
$ws.Cells.Item(9, 10).Value = 0.1154111844176443
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.08217633333333334
$ws.Cells.Item(9, 14).Value = 0.246529
$ws.Cells.Item(9, 15).Value = 0.0005826407504923545
$ws.Cells.Item(9, 16).Value = 0.0005826407504923544
$ws.Cells.Item(9, 17).Value = 0.05246846575677779
$ws.Cells.Item(9, 18).Value = 0.472216191811
$ws.Cells.Item(9, 19).Value = [double]"6.724325910430779E-05"
$ws.Cells.Item(9, 20).Value = [double]"6.724325910430777E-05"

# Row 10: MuSCs -> ECs
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Pgf"
$ws.Cells.Item(10, 3).Value = "Flt1"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 1.109135333333333
$ws.Cells.Item(10, 8).Value = 3.327406
$ws.Cells.Item(10, 9).Value = 0.2004845144158011
$ws.Cells.Item(10, 10).Value = 0.2004845144158011
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 136.544502
$ws.Cells.Item(10, 14).Value = 409.633506
$ws.Cells.Item(10, 15).Value = 0.9681180443787725
$ws.Cells.Item(10, 16).Value = 0.9681180443787725
$ws.Cells.Item(10, 17).Value = 151.446331740604
$ws.Cells.Item(10, 18).Value = 1363.016985665436
$ws.Cells.Item(10, 19).Value = 0.1940926760244532
$ws.Cells.Item(10, 20).Value = 0.1940926760244532

# Row 11: MuSCs -> FAPs
$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "Pgf"
$ws.Cells.Item(11, 3).Value = "Flt1"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 1.109135333333333
$ws.Cells.Item(11, 8).Value = 3.327406
$ws.Cells.Item(11, 9).Value = 0.2004845144158011
$ws.Cells.Item(11, 10).Value = 0.2004845144158011
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 0.3890956666666667
$ws.Cells.Item(11, 14).Value = 1.167287
$ws.Cells.Item(11, 15).Value = 0.002758738216274633
$ws.Cells.Item(11, 16).Value = 0.002758738216274633
$ws.Cells.Item(11, 17).Value = 0.4315597519468889
$ws.Cells.Item(11, 18).Value = 3.884037767522
$ws.Cells.Item(11, 19).Value = 0.0005530842916901332
$ws.Cells.Item(11, 20).Value = 0.0005530842916901331

# Row 12: MuSCs -> MuSCs
$ws.Cells.Item(12, 1).Value = "MuSCs"
$ws.Cells.Item(12, 2).Value = "Pgf"
$ws.Cells.Item(12, 3).Value = "Flt1"
$ws.Cells.Item(12, 4).Value = "MuSCs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 1.109135333333333
$ws.Cells.Item(12, 8).Value = 3.327406
$ws.Cells.Item(12, 9).Value = 0.2004845144158011
$ws.Cells.Item(12, 10).Value = 0.2004845144158011
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 4.025396333333334
$ws.Cells.Item(12, 14).Value = 12.076189
$ws.Cells.Item(12, 15).Value = 0.0285405766544606
$ws.Cells.Item(12, 16).Value = 0.02854057665446059
$ws.Cells.Item(12, 17).Value = 4.464709303970445
$ws.Cells.Item(12, 18).Value = 40.182383735734
$ws.Cells.Item(12, 19).Value = 0.005721943651716484
$ws.Cells.Item(12, 20).Value = 0.005721943651716482

# Row 13: MuSCs -> Resolving-Mac
$ws.Cells.Item(13, 1).Value = "MuSCs"
$ws.Cells.Item(13, 2).Value = "Pgf"
$ws.Cells.Item(13, 3).Value = "Flt1"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 1.109135333333333
$ws.Cells.Item(13, 8).Value = 3.327406
$ws.Cells.Item(13, 9).Value = 0.2004845144158011
$ws.Cells.Item(13, 10).Value = 0.2004845144158011
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.08217633333333334
$ws.Cells.Item(13, 14).Value = 0.246529
$ws.Cells.Item(13, 15).Value = 0.0005826407504923545
$ws.Cells.Item(13, 16).Value = 0.0005826407504923544
$ws.Cells.Item(13, 17).Value = 0.09114467486377778
$ws.Cells.Item(13, 18).Value = 0.820302073774
$ws.Cells.Item(13, 19).Value = 0.0001168104479413176
$ws.Cells.Item(13, 20).Value = 0.0001168104479413176

# Row 14: Resolving-Mac -> ECs
$ws.Cells.Item(14, 1).Value = "Resolving-Mac"
$ws.Cells.Item(14, 2).Value = "Pgf"
$ws.Cells.Item(14, 3).Value = "Flt1"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = 0.3333333333333333
$ws.Cells.Item(14, 7).Value = 0.2551043333333333
$ws.Cells.Item(14, 8).Value = 0.765313
$ws.Cells.Item(14, 9).Value = 0.04611201794463916
$ws.Cells.Item(14, 10).Value = 0.04611201794463916
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 136.544502
$ws.Cells.Item(14, 14).Value = 409.633506
$ws.Cells.Item(14, 15).Value = 0.9681180443787725
$ws.Cells.Item(14, 16).Value = 0.9681180443787725
$ws.Cells.Item(14, 17).Value = 34.833094153042
$ws.Cells.Item(14, 18).Value = 313.497847377378
$ws.Cells.Item(14, 19).Value = 0.04464187663492292
$ws.Cells.Item(14, 20).Value = 0.04464187663492293

# Row 15: Resolving-Mac -> FAPs
$ws.Cells.Item(15, 1).Value = "Resolving-Mac"
$ws.Cells.Item(15, 2).Value = "Pgf"
$ws.Cells.Item(15, 3).Value = "Flt1"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = 0.3333333333333333
$ws.Cells.Item(15, 7).Value = 0.2551043333333333
$ws.Cells.Item(15, 8).Value = 0.765313
$ws.Cells.Item(15, 9).Value = 0.04611201794463916
$ws.Cells.Item(15, 10).Value = 0.04611201794463916
$ws.Cells.Item(15, 11).Value = 2
$ws.Cells.Item(15, 12).Value = 0.6666666666666666
$ws.Cells.Item(15, 13).Value = 0.3890956666666667
$ws.Cells.Item(15, 14).Value = 1.167287
$ws.Cells.Item(15, 15).Value = 0.002758738216274633
$ws.Cells.Item(15, 16).Value = 0.002758738216274633
$ws.Cells.Item(15, 17).Value = 0.09925999064788889
$ws.Cells.Item(15, 18).Value = 0.893339915831
$ws.Cells.Item(15, 19).Value = 0.0001272109861334177
$ws.Cells.Item(15, 20).Value = 0.0001272109861334177

# Row 16: Resolving-Mac -> MuSCs
$ws.Cells.Item(16, 1).Value = "Resolving-Mac"
$ws.Cells.Item(16, 2).Value = "Pgf"
$ws.Cells.Item(16, 3).Value = "Flt1"
$ws.Cells.Item(16, 4).Value = "MuSCs"
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 6).Value = 0.3333333333333333
$ws.Cells.Item(16, 7).Value = 0.2551043333333333
$ws.Cells.Item(16, 8).Value = 0.765313
$ws.Cells.Item(16, 9).Value = 0.04611201794463916
$ws.Cells.Item(16, 10).Value = 0.04611201794463916
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 4.025396333333334
$ws.Cells.Item(16, 14).Value = 12.076189
$ws.Cells.Item(16, 15).Value = 0.0285405766544606
$ws.Cells.Item(16, 16).Value = 0.02854057665446059
$ws.Cells.Item(16, 17).Value = 1.026896048017444
$ws.Cells.Item(16, 18).Value = 9.242064432157001
$ws.Cells.Item(16, 19).Value = 0.001316063582840836
$ws.Cells.Item(16, 20).Value = 0.001316063582840836

# Row 17: Resolving-Mac -> Resolving-Mac
$ws.Cells.Item(17, 1).Value = "Resolving-Mac"
$ws.Cells.Item(17, 2).Value = "Pgf"
$ws.Cells.Item(17, 3).Value = "Flt1"
$ws.Cells.Item(17, 4).Value = "Resolving-Mac"
$ws.Cells.Item(17, 5).Value = 1
$ws.Cells.Item(17, 6).Value = 0.3333333333333333
$ws.Cells.Item(17, 7).Value = 0.2551043333333333
$ws.Cells.Item(17, 8).Value = 0.765313
$ws.Cells.Item(17, 9).Value = 0.04611201794463916
$ws.Cells.Item(17, 10).Value = 0.04611201794463916
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 0.08217633333333334
$ws.Cells.Item(17, 14).Value = 0.246529
$ws.Cells.Item(17, 15).Value = 0.0005826407504923545
$ws.Cells.Item(17, 16).Value = 0.0005826407504923544
$ws.Cells.Item(17, 17).Value = 0.02096353873077778
$ws.Cells.Item(17, 18).Value = 0.188671848577
$ws.Cells.Item(17, 19).Value = [double]"2.686674074198148E-05"
$ws.Cells.Item(17, 20).Value = [double]"2.686674074198148E-05"

